# Fix issues for modern sulfur cycle preset
#
# The NO2 column (column O) turned out to be unused / should be removed from
# the profile, so delete it outright (this also shifts every later species
# column one position to the left and removes the now-redundant NO2 entry
# from the shared-string table). The trailing "M" (total mass) column keeps
# recomputing its SUM() formula automatically against the new column range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "NO2" column (column O, the 15th column).
$ws.Columns.Item(15).Delete()

# Reflect the user's new selection: the whole of (the now former-P, now-O) column.
$ws.Range("O:O").Select() | Out-Null

# Match the print/page setup that was saved along with the workbook.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
